$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# existing Late / heading("Original") / Outstanding columns one place
# to the right (N->O, O->P, P->Q).
$ws.Columns("N").Insert()

# The newly inserted column picks up the width Excel would assign when
# inserting next to column M (10.7109375 characters); the closest
# achievable value through the exposed ColumnWidth grid is used.
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and restore the
# author's final selection on it.
$ws.Activate()
$ws.Range("C18").Select()
